$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (style) from the existing header cell H1 to the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 3

$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 4
